# Export database to Excel
# Applies updates to the C2Coverview sheet: new ECHA/readacroos/PHAROS
# reference columns (H) for each test-type block, a shift of the
# in-vivo/in-vitro OECD test-method dropdown selections (C33:C38) to make
# room for the new "In vivo mam: OECD 489" list entry, and refreshed
# window/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C2Coverview")
$ws.Activate()

# --- New shared string must be created first so it lands at the same
# --- sharedStrings index the source workbook used (index 237), ahead of
# --- ECHA / readacroos / PHAROS (238-240).
$ws.Range("C37").Value = "In vivo mam: OECD 489"

# --- Shift the OECD in-vitro/in-vivo test selections down one slot so the
# --- newly added dropdown entry has a home, pushing the former last pick
# --- ("In vivo mam: OECD 486") out to "No data".
$ws.Range("C33").Value = "In vitro mam: OECD 476"
$ws.Range("C34").Value = "In vitro mam: OECD 490"
$ws.Range("C35").Value = "In vivo mam: OECD 488"
$ws.Range("C36").Value = "In vivo mam: OECD 486"
$ws.Range("C38").Value = "No data"

$ws.Range("H18").Value = "ECHA"
$ws.Range("H19").Value = "readacroos"
$ws.Range("H20").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("H32").Value = "ECHA"
$ws.Range("H33").Value = "readacroos"
$ws.Range("H34").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H54").PasteSpecial(-4122)
$ws.Range("H54").Value = "ECHA"
$ws.Range("H55").Value = "readacroos"
$ws.Range("H56").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H62").PasteSpecial(-4122)
$ws.Range("H62").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H63").PasteSpecial(-4122)
$ws.Range("H63").Value = "readacroos"
$ws.Range("H19").Copy()
$ws.Range("H64").PasteSpecial(-4122)
$ws.Range("H64").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H70").PasteSpecial(-4122)
$ws.Range("H70").Value = "ECHA"
$ws.Range("H71").Value = "readacroos"
$ws.Range("H72").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H76").PasteSpecial(-4122)
$ws.Range("H76").Value = "ECHA"
$ws.Range("H77").Value = "readacroos"
$ws.Range("H78").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H88").PasteSpecial(-4122)
$ws.Range("H88").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H89").PasteSpecial(-4122)
$ws.Range("H89").Value = "readacroos"
$ws.Range("H90").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H96").PasteSpecial(-4122)
$ws.Range("H96").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H97").PasteSpecial(-4122)
$ws.Range("H97").Value = "readacroos"
$ws.Range("H19").Copy()
$ws.Range("H98").PasteSpecial(-4122)
$ws.Range("H98").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H104").PasteSpecial(-4122)
$ws.Range("H104").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H105").PasteSpecial(-4122)
$ws.Range("H105").Value = "readacroos"
$ws.Range("H19").Copy()
$ws.Range("H106").PasteSpecial(-4122)
$ws.Range("H106").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H117").PasteSpecial(-4122)
$ws.Range("H117").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H118").PasteSpecial(-4122)
$ws.Range("H118").Value = "readacroos"
$ws.Range("H19").Copy()
$ws.Range("H119").PasteSpecial(-4122)
$ws.Range("H119").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H125").PasteSpecial(-4122)
$ws.Range("H125").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H126").PasteSpecial(-4122)
$ws.Range("H126").Value = "readacroos"
$ws.Range("H19").Copy()
$ws.Range("H127").PasteSpecial(-4122)
$ws.Range("H127").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H130").PasteSpecial(-4122)
$ws.Range("H130").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H131").PasteSpecial(-4122)
$ws.Range("H131").Value = "readacroos"
$ws.Range("H19").Copy()
$ws.Range("H132").PasteSpecial(-4122)
$ws.Range("H132").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H143").PasteSpecial(-4122)
$ws.Range("H143").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H144").PasteSpecial(-4122)
$ws.Range("H144").Value = "readacroos"
$ws.Range("H19").Copy()
$ws.Range("H145").PasteSpecial(-4122)
$ws.Range("H145").Value = "PHAROS"
$ws.Range("H18").Copy()
$ws.Range("H149").PasteSpecial(-4122)
$ws.Range("H149").Value = "ECHA"
$ws.Range("H19").Copy()
$ws.Range("H150").PasteSpecial(-4122)
$ws.Range("H150").Value = "readacroos"
$ws.Range("H19").Copy()
$ws.Range("H151").PasteSpecial(-4122)
$ws.Range("H151").Value = "PHAROS"

# --- Refresh the active selection/view to where the edits were made.
$ws.Range("H163").Select()
